# remove Gamelogic project, modify SLG building config
# Adds two new columns (Icon, ShowName) to the BB_Build sheet between the
# existing "UpStateFunc" (F) and "Desc" (old G, now I) columns, and fills
# them in with per-row data derived from the Prefab path / existing Desc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the old column G ("Desc") out to the new column I ---------
# Read the old value first, then write it into I. Row 1 is the header
# ("Desc") and keeps the default (no) style; data rows 2-10 use the same
# "@" text style as the rest of the table.
for ($r = 1; $r -le 10; $r++) {
    $oldDesc = $ws.Cells.Item($r, 7).Value()
    if ($r -ge 2) {
        $ws.Cells.Item($r, 9).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 9).Value = $oldDesc
}

# --- 2. Header row (row 1): new G/H headers ------------------------------
$ws.Cells.Item(1, 7).Value = "Icon"
$ws.Cells.Item(1, 8).Value = "ShowName"

# --- 3. Data rows: fill in Icon (G) and ShowName (H) ---------------------
# ShowName duplicates the (already-moved) Desc text in column I.
# Icon is the prefab object's short name (last path segment of column D).
$iconNames = @{
    2  = "Altar_1_1"
    3  = "Arena_1_1"
    4  = "Camp_1_1"
    5  = "GoldMine_1_1"
    6  = "Item_hourse_1_1"
    7  = "League_1_1"
    8  = "MagicHourse_1_1"
    9  = "Tower_1_1"
    10 = "Town_1_1"
}

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $iconNames[$r]
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 9).Value()
}

# --- 4. Column widths: G:I all become width 11 ---------------------------
# (ColumnWidth=10.25 round-trips to the stored OOXML width of 11 for this
# workbook's default font/MDW, matching columns 7-9 in the target file.)
$ws.Range("G1:I1").EntireColumn.ColumnWidth = 10.25

# --- 5. Selection moves to H10 -------------------------------------------
$ws.Range("H10").Select()
